# "dev: enhance modbus byte order"
# Adds a new "weight" column (Modbus word/byte-order weight, values 1..10)
# between the existing "order" and "frequency" columns, i.e. a new column G
# is appended holding the data that used to live in column F ("frequency"),
# and column F itself becomes the new "weight" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at G (i.e. immediately to the right of F).
# This leaves column F completely untouched (content + exact width), while
# giving us an empty column G to receive the old "frequency" data.
$ws.Columns("G").Insert()

# Move the old "frequency" header + values from F into the new column G.
$ws.Range("G1").Value = $ws.Range("F1").Value()
$ws.Range("G2").Value = $ws.Range("F2").Value()
$ws.Range("G3").Value = $ws.Range("F3").Value()
$ws.Range("G4").Value = $ws.Range("F4").Value()
$ws.Range("G5").Value = $ws.Range("F5").Value()
$ws.Range("G6").Value = $ws.Range("F6").Value()
$ws.Range("G7").Value = $ws.Range("F7").Value()
$ws.Range("G8").Value = $ws.Range("F8").Value()
$ws.Range("G9").Value = $ws.Range("F9").Value()
$ws.Range("G10").Value = $ws.Range("F10").Value()
$ws.Range("G11").Value = $ws.Range("F11").Value()

# G keeps the exact same visual width as F (they're a matched pair of
# "number" columns in the new layout).
$ws.Columns("G").ColumnWidth = 13.072

# Turn column F into the new "weight" column.
$ws.Range("F1").Value = "weight"
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 4
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 6
$ws.Range("F8").Value = 7
$ws.Range("F9").Value = 8
$ws.Range("F10").Value = 9
$ws.Range("F11").Value = 10

# Column A got a little narrower in the new layout.
$ws.Columns("A").ColumnWidth = 17.501

# Match the saved selection/active cell from the edit.
$ws.Range("F2:F11").Select()
